$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E1").Value = "03_03_2024"
$ws.Range("E2").Select()
